$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cthrc1"
$ws.Range("C2").Value = "Fzd6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01838633333333333
$ws.Range("H2").Value = 0.055159
$ws.Range("I2").Value = 0.003339500866342531
$ws.Range("J2").Value = 0.003339500866342531
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.502639
$ws.Range("N2").Value = 31.507917
$ws.Range("O2").Value = 0.9701256668284471
$ws.Range("P2").Value = 0.970125666828447
$ws.Range("Q2").Value = 0.1931050215336667
$ws.Range("R2").Value = 1.737945193803
$ws.Range("S2").Value = 0.003239735504834725
$ws.Range("T2").Value = 0.003239735504834724

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cthrc1"
$ws.Range("C3").Value = "Fzd6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01838633333333333
$ws.Range("H3").Value = 0.055159
$ws.Range("I3").Value = 0.003339500866342531
$ws.Range("J3").Value = 0.003339500866342531
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2927206666666667
$ws.Range("N3").Value = 0.878162
$ws.Range("O3").Value = 0.02703852164627077
$ws.Range("P3").Value = 0.02703852164627077
$ws.Range("Q3").Value = 0.00538205975088889
$ws.Range("R3").Value = 0.048438537758
$ws.Range("S3").Value = 0.00009029516646234253
$ws.Range("T3").Value = 0.0000902951664623425

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cthrc1"
$ws.Range("C4").Value = "Fzd6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01838633333333333
$ws.Range("H4").Value = 0.055159
$ws.Range("I4").Value = 0.003339500866342531
$ws.Range("J4").Value = 0.003339500866342531
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03070066666666667
$ws.Range("N4").Value = 0.092102
$ws.Range("O4").Value = 0.002835811525282158
$ws.Range("P4").Value = 0.002835811525282158
$ws.Range("Q4").Value = 0.000564472690888889
$ws.Range("R4").Value = 0.005080254218
$ws.Range("S4").Value = 0.000009470195045463901
$ws.Range("T4").Value = 0.000009470195045463901

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cthrc1"
$ws.Range("C5").Value = "Fzd6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.449420333333333
$ws.Range("H5").Value = 16.348261
$ws.Range("I5").Value = 0.9897755900704113
$ws.Range("J5").Value = 0.9897755900704112
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.502639
$ws.Range("N5").Value = 31.507917
$ws.Range("O5").Value = 0.9701256668284471
$ws.Range("P5").Value = 0.970125666828447
$ws.Range("Q5").Value = 57.23329452025967
$ws.Range("R5").Value = 515.099650682337
$ws.Range("S5").Value = 0.9602067043275775
$ws.Range("T5").Value = 0.9602067043275773

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cthrc1"
$ws.Range("C6").Value = "Fzd6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.449420333333333
$ws.Range("H6").Value = 16.348261
$ws.Range("I6").Value = 0.9897755900704113
$ws.Range("J6").Value = 0.9897755900704112
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2927206666666667
$ws.Range("N6").Value = 0.878162
$ws.Range("O6").Value = 0.02703852164627077
$ws.Range("P6").Value = 0.02703852164627077
$ws.Range("Q6").Value = 1.595157952920222
$ws.Range("R6").Value = 14.356421576282
$ws.Range("S6").Value = 0.02676206871706924
$ws.Range("T6").Value = 0.02676206871706923

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cthrc1"
$ws.Range("C7").Value = "Fzd6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.449420333333333
$ws.Range("H7").Value = 16.348261
$ws.Range("I7").Value = 0.9897755900704113
$ws.Range("J7").Value = 0.9897755900704112
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03070066666666667
$ws.Range("N7").Value = 0.092102
$ws.Range("O7").Value = 0.002835811525282158
$ws.Range("P7").Value = 0.002835811525282158
$ws.Range("Q7").Value = 0.1673008371802222
$ws.Range("R7").Value = 1.505707534622
$ws.Range("S7").Value = 0.002806817025764621
$ws.Range("T7").Value = 0.002806817025764621

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cthrc1"
$ws.Range("C8").Value = "Fzd6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.03790633333333333
$ws.Range("H8").Value = 0.113719
$ws.Range("I8").Value = 0.006884909063246366
$ws.Range("J8").Value = 0.006884909063246365
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.502639
$ws.Range("N8").Value = 31.507917
$ws.Range("O8").Value = 0.9701256668284471
$ws.Range("P8").Value = 0.970125666828447
$ws.Range("Q8").Value = 0.3981165348136667
$ws.Range("R8").Value = 3.583048813323
$ws.Range("S8").Value = 0.006679226996035101
$ws.Range("T8").Value = 0.006679226996035098

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cthrc1"
$ws.Range("C9").Value = "Fzd6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.03790633333333333
$ws.Range("H9").Value = 0.113719
$ws.Range("I9").Value = 0.006884909063246366
$ws.Range("J9").Value = 0.006884909063246365
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2927206666666667
$ws.Range("N9").Value = 0.878162
$ws.Range("O9").Value = 0.02703852164627077
$ws.Range("P9").Value = 0.02703852164627077
$ws.Range("Q9").Value = 0.01109596716422222
$ws.Range("R9").Value = 0.099863704478
$ws.Range("S9").Value = 0.0001861577627391927
$ws.Range("T9").Value = 0.0001861577627391926

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Cthrc1"
$ws.Range("C10").Value = "Fzd6"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.03790633333333333
$ws.Range("H10").Value = 0.113719
$ws.Range("I10").Value = 0.006884909063246366
$ws.Range("J10").Value = 0.006884909063246365
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03070066666666667
$ws.Range("N10").Value = 0.092102
$ws.Range("O10").Value = 0.002835811525282158
$ws.Range("P10").Value = 0.002835811525282158
$ws.Range("Q10").Value = 0.001163749704222222
$ws.Range("R10").Value = 0.010473747338
$ws.Range("S10").Value = 0.00001952430447207363
$ws.Range("T10").Value = 0.00001952430447207363
